$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "log_elastic_mod_mean" column (W). This shifts the
# subsequent columns (cross_section, char_strength_disp, weibull_modulus_disp)
# one position to the left.
$ws.Range("W1").EntireColumn.Delete()
